# Finished off Grassland Warrior and Vengeful Arbor (Beast-Warrior and Plant searchers)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Single-Type Searchers")

# --- Row 4: Beast-Warrior searcher -> "Grassland Warrior" ---
# Copy the formatting of an already-finished row (row 7) onto row 4,
# then fill in the new card's text/ATK/DEF values.
$ws.Range("E7:L7").Copy() | Out-Null
$ws.Range("E4:L4").PasteSpecial(-4122) | Out-Null
$ws.Range("E4").Value = "Grassland Warrior"
$ws.Range("H4").Value = 900
$ws.Range("I4").Value = 1500

# --- Row 12: Plant searcher -> "Vengeful Arbor" ---
$ws.Range("E7:L7").Copy() | Out-Null
$ws.Range("E12:L12").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").Value = "Vengeful Arbor"
$ws.Range("H12").Value = 800
$ws.Range("I12").Value = 1200

$excel.CutCopyMode = 0

# --- Selection reflects the just-finished card (row 12) ---
$ws.Range("E12:L12").Select() | Out-Null
